# Updated cryptos list values (Price and Volume(1h) columns)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '56.945.04'
$ws.Range("E2").Value = '  -0.71%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.318.37'

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.01'
$ws.Range("E4").Value = '  +0.29%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '530.79'
$ws.Range("E5").Value = '  +2.53%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '132.70'
$ws.Range("E6").Value = '  -2.13%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.996'
$ws.Range("E7").Value = '  -0.22%  '

$ws.Range("E8").Value = '  -0.95%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.342.19'
$ws.Range("E9").Value = '  -1.91%  '

$ws.Range("E10").Value = '  -1.32%  '

$ws.Range("E11").Value = '  +0.21%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.39'
$ws.Range("E12").Value = '  -2.84%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.352'
$ws.Range("E13").Value = '  +2.54%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.736.67'
$ws.Range("E14").Value = '  -2.45%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '23.32'
$ws.Range("E15").Value = '  -4.17%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '57.276.58'
$ws.Range("E16").Value = '  -0.23%  '

$ws.Range("E17").Value = '  -2.10%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.338.11'
$ws.Range("E18").Value = '  -1.46%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '337.55'
$ws.Range("E19").Value = '  +2.40%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.41'
$ws.Range("E20").Value = '  -1.68%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.19'
$ws.Range("E21").Value = '  -1.59%  '

$ws.Range("E22").Value = '  +0.82%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.998'
$ws.Range("E23").Value = '  -0.05%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '62.09'
$ws.Range("E24").Value = '  +1.10%  '

$ws.Range("E25").Value = '  +0.63%  '

$ws.Range("E26").Value = '  -2.79%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.995'
$ws.Range("E27").Value = '  -0.08%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.36'
$ws.Range("E28").Value = '  +1.38%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '173.16'
$ws.Range("E29").Value = '  +3.60%  '

$ws.Range("E30").Value = '  +1.40%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0₃0725'
$ws.Range("E31").Value = '  -2.98%  '

$ws.Range("E32").Value = '  -2.45%  '

$ws.Range("E33").Value = '  -0.62%  '

$ws.Range("E34").Value = '  -0.04%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.992'
$ws.Range("E35").Value = '  -0.20%  '

$ws.Range("E36").Value = '  -3.69%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.933'
$ws.Range("E37").Value = '  +1.71%  '

$ws.Range("E38").Value = '  -1.65%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '39.22'
$ws.Range("E39").Value = '  +0.87%  '

$ws.Range("E40").Value = '  -1.69%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.95'
$ws.Range("E41").Value = '  +11.69%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '149.45'
$ws.Range("E42").Value = '  -0.80%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.374'
$ws.Range("E43").Value = '  -3.38%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.61'
$ws.Range("E44").Value = '  -1.56%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '284.52'
$ws.Range("E45").Value = '  -1.13%  '

$ws.Range("E46").Value = '  -1.22%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '18.75'
$ws.Range("E48").Value = '  +3.01%  '

$ws.Range("E49").Value = '  -1.63%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0217'
$ws.Range("E50").Value = '  -1.09%  '

$ws.Range("E51").Value = '  -1.84%  '
